# Update NATMI TPM-derived ligand/receptor edge-weight values (Spp1-Itgav)
# per the new TPM recomputation. Source columns G,H,I,J (ligand) and
# M,N,O,P (receptor) are raw inputs; Q,R,S,T (edge weights/specificity)
# are the corresponding products (G*M, H*N, I*O, J*P) from the pipeline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 20.453651
$ws.Range("H2").Value = 40.907302
$ws.Range("I2").Value = 0.002931177245748417
$ws.Range("J2").Value = 0.001974175492708352
$ws.Range("M2").Value = 20.574342
$ws.Range("N2").Value = 41.148684
$ws.Range("O2").Value = 0.07442291871210138
$ws.Range("P2").Value = 0.05295769307665528
$ws.Range("Q2").Value = 420.820410822642
$ws.Range("R2").Value = 1683.281643290568
$ws.Range("S2").Value = 0.0002181467658910957
$ws.Range("T2").Value = 0.0001045477798223036

# Row 3
$ws.Range("G3").Value = 20.453651
$ws.Range("H3").Value = 40.907302
$ws.Range("I3").Value = 0.002931177245748417
$ws.Range("J3").Value = 0.001974175492708352
$ws.Range("O3").Value = 0.2107257583291067
$ws.Range("P3").Value = 0.2249216416002446
$ws.Range("Q3").Value = 1191.537522654943
$ws.Range("R3").Value = 7149.225135929659
$ws.Range("S3").Value = 0.0006176745479073576
$ws.Range("T3").Value = 0.0004440347926269344

# Row 4
$ws.Range("G4").Value = 20.453651
$ws.Range("H4").Value = 40.907302
$ws.Range("I4").Value = 0.002931177245748417
$ws.Range("J4").Value = 0.001974175492708352
$ws.Range("M4").Value = 59.53576899999999
$ws.Range("N4").Value = 178.607307
$ws.Range("O4").Value = 0.2153568603433074
$ws.Range("P4").Value = 0.2298647253300724
$ws.Range("Q4").Value = 1217.723841142619
$ws.Range("R4").Value = 7306.343046855714
$ws.Range("S4").Value = 0.0006312491287541224
$ws.Range("T4").Value = 0.0004537933073847657

# Row 5
$ws.Range("G5").Value = 20.453651
$ws.Range("H5").Value = 40.907302
$ws.Range("I5").Value = 0.002931177245748417
$ws.Range("J5").Value = 0.001974175492708352
$ws.Range("M5").Value = 31.770234
$ws.Range("N5").Value = 63.540468
$ws.Range("O5").Value = 0.1149214658940947
$ws.Range("P5").Value = 0.08177555817559162
$ws.Range("Q5").Value = 649.8172784243341
$ws.Range("R5").Value = 2599.269113697336
$ws.Range("S5").Value = 0.000336855185876823
$ws.Range("T5").Value = 0.0001614393028527991

# Row 6
$ws.Range("G6").Value = 20.453651
$ws.Range("H6").Value = 40.907302
$ws.Range("I6").Value = 0.002931177245748417
$ws.Range("J6").Value = 0.001974175492708352
$ws.Range("M6").Value = 31.46548433333334
$ws.Range("N6").Value = 94.39645300000001
$ws.Range("O6").Value = 0.1138191045320704
$ws.Range("P6").Value = 0.121486713536183
$ws.Range("Q6").Value = 643.5840350999678
$ws.Range("R6").Value = 3861.504210599806
$ws.Range("S6").Value = 0.0003336239693358653
$ws.Range("T6").Value = 0.0002398360925528124

# Row 7
$ws.Range("G7").Value = 20.453651
$ws.Range("H7").Value = 40.907302
$ws.Range("I7").Value = 0.002931177245748417
$ws.Range("J7").Value = 0.001974175492708352
$ws.Range("M7").Value = 74.85037233333333
$ws.Range("N7").Value = 224.551117
$ws.Range("O7").Value = 0.2707538921893195
$ws.Range("P7").Value = 0.2889936682812531
$ws.Range("Q7").Value = 1530.963392926056
$ws.Range("R7").Value = 9185.780357556334
$ws.Range("S7").Value = 0.0007936276479831534
$ws.Range("T7").Value = 0.000570524217468737

# Row 8
$ws.Range("I8").Value = 0.004554291484292619
$ws.Range("J8").Value = 0.004601037338145601
$ws.Range("M8").Value = 20.574342
$ws.Range("N8").Value = 41.148684
$ws.Range("O8").Value = 0.07442291871210138
$ws.Range("P8").Value = 0.05295769307665528
$ws.Range("Q8").Value = 653.8461009841561
$ws.Range("R8").Value = 3923.076605904937
$ws.Range("S8").Value = 0.0003389436649267251
$ws.Range("T8").Value = 0.0002436603231877457

# Row 9
$ws.Range("I9").Value = 0.004554291484292619
$ws.Range("J9").Value = 0.004601037338145601
$ws.Range("O9").Value = 0.2107257583291067
$ws.Range("P9").Value = 0.2249216416002446
$ws.Range("S9").Value = 0.0009597065266793552
$ws.Range("T9").Value = 0.001034872871159728

# Row 10
$ws.Range("I10").Value = 0.004554291484292619
$ws.Range("J10").Value = 0.004601037338145601
$ws.Range("M10").Value = 59.53576899999999
$ws.Range("N10").Value = 178.607307
$ws.Range("O10").Value = 0.2153568603433074
$ws.Range("P10").Value = 0.2298647253300724
$ws.Range("Q10").Value = 1892.027965207508
$ws.Range("R10").Value = 17028.25168686758
$ws.Range("S10").Value = 0.00098079791514552
$ws.Range("T10").Value = 0.001057616183966246

# Row 11
$ws.Range("I11").Value = 0.004554291484292619
$ws.Range("J11").Value = 0.004601037338145601
$ws.Range("M11").Value = 31.770234
$ws.Range("N11").Value = 63.540468
$ws.Range("O11").Value = 0.1149214658940947
$ws.Range("P11").Value = 0.08177555817559162
$ws.Range("Q11").Value = 1009.648018306212
$ws.Range("R11").Value = 6057.888109837273
$ws.Range("S11").Value = 0.0005233858534839
$ws.Range("T11").Value = 0.0003762523965135948

# Row 12
$ws.Range("I12").Value = 0.004554291484292619
$ws.Range("J12").Value = 0.004601037338145601
$ws.Range("M12").Value = 31.46548433333334
$ws.Range("N12").Value = 94.39645300000001
$ws.Range("O12").Value = 0.1138191045320704
$ws.Range("P12").Value = 0.121486713536183
$ws.Range("Q12").Value = 999.9631699972737
$ws.Range("R12").Value = 8999.668529975463
$ws.Range("S12").Value = 0.0005183653785202196
$ws.Range("T12").Value = 0.0005589649050685763

# Row 13
$ws.Range("I13").Value = 0.004554291484292619
$ws.Range("J13").Value = 0.004601037338145601
$ws.Range("M13").Value = 74.85037233333333
$ws.Range("N13").Value = 224.551117
$ws.Range("O13").Value = 0.2707538921893195
$ws.Range("P13").Value = 0.2889936682812531
$ws.Range("Q13").Value = 2378.721229935924
$ws.Range("R13").Value = 21408.49106942332
$ws.Range("S13").Value = 0.0012330921455369
$ws.Range("T13").Value = 0.00132967065824971

# Row 14
$ws.Range("G14").Value = 917.50354
$ws.Range("H14").Value = 2752.51062
$ws.Range("I14").Value = 0.131485840808647
$ws.Range("J14").Value = 0.1328354289760658
$ws.Range("M14").Value = 20.574342
$ws.Range("N14").Value = 41.148684
$ws.Range("O14").Value = 0.07442291871210138
$ws.Range("P14").Value = 0.05295769307665528
$ws.Range("Q14").Value = 18877.03161817068
$ws.Range("R14").Value = 113262.1897090241
$ws.Range("S14").Value = 0.009785560042294236
$ws.Range("T14").Value = 0.007034657877420332

# Row 15
$ws.Range("G15").Value = 917.50354
$ws.Range("H15").Value = 2752.51062
$ws.Range("I15").Value = 0.131485840808647
$ws.Range("J15").Value = 0.1328354289760658
$ws.Range("O15").Value = 0.2107257583291067
$ws.Range("P15").Value = 0.2249216416002446
$ws.Range("Q15").Value = 53449.62105194522
$ws.Range("R15").Value = 481046.589467507
$ws.Range("S15").Value = 0.02770745351394234
$ws.Range("T15").Value = 0.02987756274796941

# Row 16
$ws.Range("G16").Value = 917.50354
$ws.Range("H16").Value = 2752.51062
$ws.Range("I16").Value = 0.131485840808647
$ws.Range("J16").Value = 0.1328354289760658
$ws.Range("M16").Value = 59.53576899999999
$ws.Range("N16").Value = 178.607307
$ws.Range("O16").Value = 0.2153568603433074
$ws.Range("P16").Value = 0.2298647253300724
$ws.Range("Q16").Value = 54624.27881412226
$ws.Range("R16").Value = 491618.5093271003
$ws.Range("S16").Value = 0.02831637785615014
$ws.Range("T16").Value = 0.03053417939568569

# Row 17
$ws.Range("G17").Value = 917.50354
$ws.Range("H17").Value = 2752.51062
$ws.Range("I17").Value = 0.131485840808647
$ws.Range("J17").Value = 0.1328354289760658
$ws.Range("M17").Value = 31.770234
$ws.Range("N17").Value = 63.540468
$ws.Range("O17").Value = 0.1149214658940947
$ws.Range("P17").Value = 0.08177555817559162
$ws.Range("Q17").Value = 29149.30216162836
$ws.Range("R17").Value = 174895.8129697702
$ws.Range("S17").Value = 0.01511054557004728
$ws.Range("T17").Value = 0.01086269135001194

# Row 18
$ws.Range("G18").Value = 917.50354
$ws.Range("H18").Value = 2752.51062
$ws.Range("I18").Value = 0.131485840808647
$ws.Range("J18").Value = 0.1328354289760658
$ws.Range("M18").Value = 31.46548433333334
$ws.Range("N18").Value = 94.39645300000001
$ws.Range("O18").Value = 0.1138191045320704
$ws.Range("P18").Value = 0.121486713536183
$ws.Range("Q18").Value = 28869.69326364788
$ws.Range("R18").Value = 259827.2393728309
$ws.Range("S18").Value = 0.01496560065948656
$ws.Range("T18").Value = 0.01613773970747128

# Row 19
$ws.Range("G19").Value = 917.50354
$ws.Range("H19").Value = 2752.51062
$ws.Range("I19").Value = 0.131485840808647
$ws.Range("J19").Value = 0.1328354289760658
$ws.Range("M19").Value = 74.85037233333333
$ws.Range("N19").Value = 224.551117
$ws.Range("O19").Value = 0.2707538921893195
$ws.Range("P19").Value = 0.2889936682812531
$ws.Range("Q19").Value = 68675.48158615139
$ws.Range("R19").Value = 618079.3342753624
$ws.Range("S19").Value = 0.03560030316672643
$ws.Range("T19").Value = 0.03838859789750711

# Row 20
$ws.Range("G20").Value = 192.2315905
$ws.Range("H20").Value = 384.463181
$ws.Range("I20").Value = 0.02754837578814793
$ws.Range("J20").Value = 0.01855409065547504
$ws.Range("M20").Value = 20.574342
$ws.Range("N20").Value = 41.148684
$ws.Range("O20").Value = 0.07442291871210138
$ws.Range("P20").Value = 0.05295769307665528
$ws.Range("Q20").Value = 3955.038486150951
$ws.Range("R20").Value = 15820.15394460381
$ws.Range("S20").Value = 0.002050230531931755
$ws.Range("T20").Value = 0.0009825818382490849

# Row 21
$ws.Range("G21").Value = 192.2315905
$ws.Range("H21").Value = 384.463181
$ws.Range("I21").Value = 0.02754837578814793
$ws.Range("J21").Value = 0.01855409065547504
$ws.Range("O21").Value = 0.2107257583291067
$ws.Range("P21").Value = 0.2249216416002446
$ws.Range("Q21").Value = 11198.54607475162
$ws.Range("R21").Value = 67191.2764485097
$ws.Range("S21").Value = 0.005805152378692675
$ws.Range("T21").Value = 0.004173216528629205

# Row 22
$ws.Range("G22").Value = 192.2315905
$ws.Range("H22").Value = 384.463181
$ws.Range("I22").Value = 0.02754837578814793
$ws.Range("J22").Value = 0.01855409065547504
$ws.Range("M22").Value = 59.53576899999999
$ws.Range("N22").Value = 178.607307
$ws.Range("O22").Value = 0.2153568603433074
$ws.Range("P22").Value = 0.2298647253300724
$ws.Range("Q22").Value = 11444.65556651059
$ws.Range("R22").Value = 68667.93339906356
$ws.Range("S22").Value = 0.005932731717293125
$ws.Range("T22").Value = 0.004264930952270032

# Row 23
$ws.Range("G23").Value = 192.2315905
$ws.Range("H23").Value = 384.463181
$ws.Range("I23").Value = 0.02754837578814793
$ws.Range("J23").Value = 0.01855409065547504
$ws.Range("M23").Value = 31.770234
$ws.Range("N23").Value = 63.540468
$ws.Range("O23").Value = 0.1149214658940947
$ws.Range("P23").Value = 0.08177555817559162
$ws.Range("Q23").Value = 6107.242612377177
$ws.Range("R23").Value = 24428.97044950871
$ws.Range("S23").Value = 0.003165899728575345
$ws.Range("T23").Value = 0.001517271119792

# Row 24
$ws.Range("G24").Value = 192.2315905
$ws.Range("H24").Value = 384.463181
$ws.Range("I24").Value = 0.02754837578814793
$ws.Range("J24").Value = 0.01855409065547504
$ws.Range("M24").Value = 31.46548433333334
$ws.Range("N24").Value = 94.39645300000001
$ws.Range("O24").Value = 0.1138191045320704
$ws.Range("P24").Value = 0.121486713536183
$ws.Range("Q24").Value = 6048.6600992495
$ws.Range("R24").Value = 36291.960595497
$ws.Range("S24").Value = 0.003135531463519966
$ws.Range("T24").Value = 0.002254075496386065

# Row 25
$ws.Range("G25").Value = 192.2315905
$ws.Range("H25").Value = 384.463181
$ws.Range("I25").Value = 0.02754837578814793
$ws.Range("J25").Value = 0.01855409065547504
$ws.Range("M25").Value = 74.85037233333333
$ws.Range("N25").Value = 224.551117
$ws.Range("O25").Value = 0.2707538921893195
$ws.Range("P25").Value = 0.2889936682812531
$ws.Range("Q25").Value = 14388.60612315386
$ws.Range("R25").Value = 86331.63673892317
$ws.Range("S25").Value = 0.007458829968135064
$ws.Range("T25").Value = 0.005362014720148653

# Row 26
$ws.Range("G26").Value = 749.9081419999999
$ws.Range("H26").Value = 2249.724426
$ws.Range("I26").Value = 0.1074680350335435
$ws.Range("J26").Value = 0.108571101246339
$ws.Range("M26").Value = 20.574342
$ws.Range("N26").Value = 41.148684
$ws.Range("O26").Value = 0.07442291871210138
$ws.Range("P26").Value = 0.05295769307665528
$ws.Range("Q26").Value = 15428.86658209256
$ws.Range("R26").Value = 92573.19949255539
$ws.Range("S26").Value = 0.007998084835450675
$ws.Range("T26").Value = 0.005749675056798085

# Row 27
$ws.Range("G27").Value = 749.9081419999999
$ws.Range("H27").Value = 2249.724426
$ws.Range("I27").Value = 0.1074680350335435
$ws.Range("J27").Value = 0.108571101246339
$ws.Range("O27").Value = 0.2107257583291067
$ws.Range("P27").Value = 0.2249216416002446
$ws.Range("Q27").Value = 43686.268516924
$ws.Range("R27").Value = 393176.416652316
$ws.Range("S27").Value = 0.02264628317858247
$ws.Range("T27").Value = 0.02441999032267292

# Row 28
$ws.Range("G28").Value = 749.9081419999999
$ws.Range("H28").Value = 2249.724426
$ws.Range("I28").Value = 0.1074680350335435
$ws.Range("J28").Value = 0.108571101246339
$ws.Range("M28").Value = 59.53576899999999
$ws.Range("N28").Value = 178.607307
$ws.Range("O28").Value = 0.2153568603433074
$ws.Range("P28").Value = 0.2298647253300724
$ws.Range("Q28").Value = 44646.35791333119
$ws.Range("R28").Value = 401817.2212199807
$ws.Range("S28").Value = 0.0231439786120885
$ws.Range("T28").Value = 0.02495666636677318

# Row 29
$ws.Range("G29").Value = 749.9081419999999
$ws.Range("H29").Value = 2249.724426
$ws.Range("I29").Value = 0.1074680350335435
$ws.Range("J29").Value = 0.108571101246339
$ws.Range("M29").Value = 31.770234
$ws.Range("N29").Value = 63.540468
$ws.Range("O29").Value = 0.1149214658940947
$ws.Range("P29").Value = 0.08177555817559162
$ws.Range("Q29").Value = 23824.75714984523
$ws.Range("R29").Value = 142948.5428990714
$ws.Range("S29").Value = 0.01235038412281274
$ws.Range("T29").Value = 0.00887846240615804

# Row 30
$ws.Range("G30").Value = 749.9081419999999
$ws.Range("H30").Value = 2249.724426
$ws.Range("I30").Value = 0.1074680350335435
$ws.Range("J30").Value = 0.108571101246339
$ws.Range("M30").Value = 31.46548433333334
$ws.Range("N30").Value = 94.39645300000001
$ws.Range("O30").Value = 0.1138191045320704
$ws.Range("P30").Value = 0.121486713536183
$ws.Range("Q30").Value = 23596.22289354011
$ws.Range("R30").Value = 212366.006041861
$ws.Range("S30").Value = 0.01223191551333909
$ws.Range("T30").Value = 0.0131899462754219

# Row 31
$ws.Range("G31").Value = 749.9081419999999
$ws.Range("H31").Value = 2249.724426
$ws.Range("I31").Value = 0.1074680350335435
$ws.Range("J31").Value = 0.108571101246339
$ws.Range("M31").Value = 74.85037233333333
$ws.Range("N31").Value = 224.551117
$ws.Range("O31").Value = 0.2707538921893195
$ws.Range("P31").Value = 0.2889936682812531
$ws.Range("Q31").Value = 56130.90364449819
$ws.Range("R31").Value = 505178.1328004837
$ws.Range("S31").Value = 0.02909738877127006
$ws.Range("T31").Value = 0.03137636081851483

# Row 32
$ws.Range("G32").Value = 5066.087972333334
$ws.Range("H32").Value = 15198.263917
$ws.Range("I32").Value = 0.7260122796396206
$ws.Range("J32").Value = 0.7334641662912662
$ws.Range("M32").Value = 20.574342
$ws.Range("N32").Value = 41.148684
$ws.Range("O32").Value = 0.07442291871210138
$ws.Range("P32").Value = 0.05295769307665528
$ws.Range("Q32").Value = 104231.4265448725
$ws.Range("R32").Value = 625388.5592692352
$ws.Range("S32").Value = 0.0540319528716069
$ws.Range("T32").Value = 0.03884257020117773

# Row 33
$ws.Range("G33").Value = 5066.087972333334
$ws.Range("H33").Value = 15198.263917
$ws.Range("I33").Value = 0.7260122796396206
$ws.Range("J33").Value = 0.7334641662912662
$ws.Range("O33").Value = 0.2107257583291067
$ws.Range("P33").Value = 0.2249216416002446
$ws.Range("Q33").Value = 295127.4524096487
$ws.Range("R33").Value = 2656147.071686838
$ws.Range("S33").Value = 0.1529894881833025
$ws.Range("T33").Value = 0.1649719643371864

# Row 34
$ws.Range("G34").Value = 5066.087972333334
$ws.Range("H34").Value = 15198.263917
$ws.Range("I34").Value = 0.7260122796396206
$ws.Range("J34").Value = 0.7334641662912662
$ws.Range("M34").Value = 59.53576899999999
$ws.Range("N34").Value = 178.607307
$ws.Range("O34").Value = 0.2153568603433074
$ws.Range("P34").Value = 0.2298647253300724
$ws.Range("Q34").Value = 301613.4432545157
$ws.Range("R34").Value = 2714520.989290642
$ws.Range("S34").Value = 0.156351725113876
$ws.Range("T34").Value = 0.1685975391239924

# Row 35
$ws.Range("G35").Value = 5066.087972333334
$ws.Range("H35").Value = 15198.263917
$ws.Range("I35").Value = 0.7260122796396206
$ws.Range("J35").Value = 0.7334641662912662
$ws.Range("M35").Value = 31.770234
$ws.Range("N35").Value = 63.540468
$ws.Range("O35").Value = 0.1149214658940947
$ws.Range("P35").Value = 0.08177555817559162
$ws.Range("Q35").Value = 160950.8003456155
$ws.Range("R35").Value = 965704.8020736933
$ws.Range("S35").Value = 0.08343439543329857
$ws.Range("T35").Value = 0.05997944160026325

# Row 36
$ws.Range("G36").Value = 5066.087972333334
$ws.Range("H36").Value = 15198.263917
$ws.Range("I36").Value = 0.7260122796396206
$ws.Range("J36").Value = 0.7334641662912662
$ws.Range("M36").Value = 31.46548433333334
$ws.Range("N36").Value = 94.39645300000001
$ws.Range("O36").Value = 0.1138191045320704
$ws.Range("P36").Value = 0.121486713536183
$ws.Range("Q36").Value = 159406.911724743
$ws.Range("R36").Value = 1434662.205522686
$ws.Range("S36").Value = 0.0826340675478687
$ws.Range("T36").Value = 0.08910615105928231

# Row 37
$ws.Range("G37").Value = 5066.087972333334
$ws.Range("H37").Value = 15198.263917
$ws.Range("I37").Value = 0.7260122796396206
$ws.Range("J37").Value = 0.7334641662912662
$ws.Range("M37").Value = 74.85037233333333
$ws.Range("N37").Value = 224.551117
$ws.Range("O37").Value = 0.2707538921893195
$ws.Range("P37").Value = 0.2889936682812531
$ws.Range("Q37").Value = 379198.5710025717
$ws.Range("R37").Value = 3412787.139023145
$ws.Range("S37").Value = 0.1965706504896679
$ws.Range("T37").Value = 0.2119664999693641
